$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the formula-driven "translation" (x1000) cells with their
#     plain numeric equivalents (no more *1000 multiplication formulas) ---
$ws.Range("C7").Value = 1.70079118954
$ws.Range("D7").Value = 3.412
$ws.Range("C8").Value = 0.0159456324149
$ws.Range("C9").Value = 1.51095763913
$ws.Range("D9").Value = 0.5

# --- Apply a 3-decimal-place number format to the rotation/translation
#     block (C3:D9) ---
$ws.Range("C3:D9").NumberFormat = "0.000_ "

# --- Move the active selection to D11, matching the saved cursor position ---
$ws.Range("D11").Select() | Out-Null
